$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $addr, $val)
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-CellText $ws 'D2' '27.501.10'
Set-CellText $ws 'E2' '  +2.81%  '
Set-CellText $ws 'D3' '1.816.77'
Set-CellText $ws 'D4' '1.003'
Set-CellText $ws 'E4' '  +0.51%  '
Set-CellText $ws 'D5' '344.15'
Set-CellText $ws 'E5' '  +3.32%  '
Set-CellText $ws 'D6' '1.000'
Set-CellText $ws 'E6' '  +0.66%  '
Set-CellText $ws 'D7' '0.3836'
Set-CellText $ws 'E7' '  +2.25%  '
Set-CellText $ws 'D8' '0.3553'
Set-CellText $ws 'E8' '  +2.89%  '
Set-CellText $ws 'D9' '48.97'
Set-CellText $ws 'E9' '  -1.96%  '
Set-CellText $ws 'D10' '1.238'
Set-CellText $ws 'E10' '  +2.41%  '
Set-CellText $ws 'D11' '0.07790'
Set-CellText $ws 'E11' '  +3.43%  '
Set-CellText $ws 'D12' '1.001'
Set-CellText $ws 'E12' '  +0.54%  '
Set-CellText $ws 'D13' '22.45'
Set-CellText $ws 'E13' '  +8.84%  '
Set-CellText $ws 'D14' '6.620'
Set-CellText $ws 'E14' '  +1.80%  '
Set-CellText $ws 'D15' '1.814.36'
Set-CellText $ws 'E15' '  +4.19%  '
Set-CellText $ws 'D16' '7.225'
Set-CellText $ws 'E16' '  +1.71%  '
Set-CellText $ws 'D17' '0.00001128'
Set-CellText $ws 'E17' '  +2.45%  '
Set-CellText $ws 'D18' '0.06737'
Set-CellText $ws 'E18' '  +0.42%  '
Set-CellText $ws 'D19' '87.04'
Set-CellText $ws 'E19' '  +3.14%  '
Set-CellText $ws 'E20' '  +0.64%  '
Set-CellText $ws 'D21' '17.67'
Set-CellText $ws 'E21' '  +4.78%  '
Set-CellText $ws 'D22' '6.571'
Set-CellText $ws 'E22' '  +5.64%  '
Set-CellText $ws 'D23' '13.18'
Set-CellText $ws 'E23' '  -0.10%  '
Set-CellText $ws 'D24' '27.494.50'
Set-CellText $ws 'E24' '  +2.87%  '
Set-CellText $ws 'D25' '2.467'
Set-CellText $ws 'E25' '  -0.35%  '
Set-CellText $ws 'D26' '2.694'
Set-CellText $ws 'E26' '  +5.89%  '
Set-CellText $ws 'E27' '  +12.38%  '
Set-CellText $ws 'D28' '1.457'
Set-CellText $ws 'E28' '  +2.80%  '
Set-CellText $ws 'D29' '154.05'
Set-CellText $ws 'E29' '  +0.55%  '
Set-CellText $ws 'D30' '2.019.95'
Set-CellText $ws 'E30' '  +4.32%  '
Set-CellText $ws 'D31' '136.16'
Set-CellText $ws 'E31' '  +2.53%  '
Set-CellText $ws 'D32' '6.409'
Set-CellText $ws 'E32' '  +2.28%  '
Set-CellText $ws 'D33' '4.098'
Set-CellText $ws 'E33' '  -0.88%  '
Set-CellText $ws 'D34' '13.98'
Set-CellText $ws 'E34' '  +5.69%  '
Set-CellText $ws 'D35' '0.08820'
Set-CellText $ws 'E35' '  +2.05%  '
Set-CellText $ws 'D36' '1.691'
Set-CellText $ws 'E36' '  -2.01%  '
Set-CellText $ws 'D37' '5.658'
Set-CellText $ws 'E37' '  +2.73%  '
Set-CellText $ws 'D38' '0.7033'
Set-CellText $ws 'E38' '  +11.43%  '
Set-CellText $ws 'B39' 'FraxShare'
Set-CellText $ws 'C39' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-CellText $ws 'D39' '9.078'
Set-CellText $ws 'E39' '  +4.24%  '
Set-CellText $ws 'B40' 'VeChain'
Set-CellText $ws 'C40' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText $ws 'D40' '0.02418'
Set-CellText $ws 'E40' '  +1.62%  '
Set-CellText $ws 'B41' 'Hedera'
Set-CellText $ws 'C41' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-CellText $ws 'D41' '0.06519'
Set-CellText $ws 'E41' '  +1.92%  '
Set-CellText $ws 'B42' 'Algorand'
Set-CellText $ws 'C42' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-CellText $ws 'D42' '0.2260'
Set-CellText $ws 'E42' '  +3.27%  '
Set-CellText $ws 'D43' '1.300'
Set-CellText $ws 'E43' '  +4.26%  '
Set-CellText $ws 'D44' '14.80'
Set-CellText $ws 'E44' '  +1.67%  '
Set-CellText $ws 'D45' '0.6623'
Set-CellText $ws 'E45' '  +8.10%  '
Set-CellText $ws 'D46' '1.000'
Set-CellText $ws 'E46' '  +0.54%  '
Set-CellText $ws 'D47' '3.966'
Set-CellText $ws 'E47' '  +1.49%  '
Set-CellText $ws 'D48' '2.195'
Set-CellText $ws 'E48' '  +5.28%  '
Set-CellText $ws 'D49' '132.78'
Set-CellText $ws 'E49' '  +2.62%  '
Set-CellText $ws 'D50' '0.07334'
Set-CellText $ws 'E50' '  -0.55%  '
Set-CellText $ws 'D51' '81.08'
Set-CellText $ws 'E51' '  +3.68%  '
